$wb = $excel.ActiveWorkbook

# Sheet 1 = "展览" and Sheet 4 = "全部类型" both contain the same event table.
# Update the "想去人数" (want-to-go count) column F for rows 2-5 on both sheets.
$sheetIndexes = @(1, 4)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Range("F2").Value = 55
    $ws.Range("F3").Value = 289
    $ws.Range("F4").Value = 22
    $ws.Range("F5").Value = 85
}
